# Update profit files after running on 2025-09-03
# Appends a new row (A17/B17) to the Date/Profit table on Sheet1, mirroring
# the existing rows: the date is stored as literal text (not an Excel date
# serial) and the profit figure as a plain number.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Write the date as a quoted literal first so Excel's input parser doesn't
# auto-convert the "09/03/2025"-looking text into a date serial value, then
# drop back to the workbook's default (Normal) style so no quote-prefix /
# number-format styling lingers on the cell - matching how the other date
# cells in this column are stored (plain text, default style).
$ws.Range("A17").Value = "'09/03/2025"
$ws.Range("A17").Style = "Normal"

$ws.Range("B17").Value = 13755.16
